$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.812.99"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "2.295.65"
$ws.Range("E3").Value = "  -5.28%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.86"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.62%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("D9").Value = "2.294.37"
$ws.Range("E9").Value = "  -5.24%  "
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.56"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.77%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.336"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.92"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "2.701.65"
$ws.Range("E15").Value = "  -5.30%  "
$ws.Range("D16").Value = "58.772.87"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.45%  "
$ws.Range("D18").Value = "2.314.05"
$ws.Range("E18").Value = "  -4.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.68"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.57"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.48"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.12%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.15"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.173"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.10"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.44%  "
$ws.Range("E28").Value = "  -7.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.75"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.62"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "0.0₃0730"
$ws.Range("E31").Value = "  -5.80%  "
$ws.Range("B32").Value = "SuiNetwork"
$ws.Range("C32").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.10"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.09%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.80"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.12%  "
$ws.Range("E34").Value = "  -4.63%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.78"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.21%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  -6.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.97"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -6.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.52"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "301.66"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -7.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "140.01"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.46"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0953"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0500"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.23%  "
$ws.Range("E47").Value = "  -3.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.63"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0216"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.69"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.02"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.18%  "
